# Re-style the three summary tables (slides 14-16) from the default
# "Table_0" style to the built-in table style
# {E9511814-E422-4334-8146-76ABD6C0F998}.
#
# PowerPoint's object model does not allow `Table.Style` to be assigned
# a GUID string directly (that raises "Table styles cannot be assigned
# through a property"); the supported call is `Table.ApplyStyle(styleId)`.

$oldStyleId = "{29C6B79C-C1B3-439A-A80F-0A475276058A}"
$newStyleId = "{E9511814-E422-4334-8146-76ABD6C0F998}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
